# Update the "想去人数" (want-to-go count) figures in column F for the
# sheets that list exhibition-type events ("展览" and "全部类型").
# Rows refer to the same underlying events in both sheets, so the same
# set of updates is applied twice.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 6540
    5  = 47
    6  = 1959
    7  = 1495
    8  = 304
    9  = 1003
    10 = 372
    11 = 6
    12 = 5623
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
